# Generate Report for Handoff
#
# Updates the localization-status report: flips the "In Translation"
# status cells to "Ready for handoff" and refreshes the associated
# generate/handback timestamps on each sheet. Excel's own column
# auto-sizing widens the status columns to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-05 00:44:54"

# Columns E:F hold the status text that just grew longer - widen them
# to match Excel's auto-fit result.
$wsOverview.Columns.Item(5).ColumnWidth = 16.38265482584637
$wsOverview.Columns.Item(6).ColumnWidth = 16.38265482584637

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-05 00:44:49"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.38265482584637

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-05 00:44:54"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.38265482584637

Write-Output "Handoff report regenerated."
